$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update usuario names (column A) and rol values (column C) for rows 2-5
$ws.Range("A2").Value = "Laura"
$ws.Range("C2").Value = "VENDEDOR"

$ws.Range("A3").Value = "Juanita"

$ws.Range("A4").Value = "Jorge"

$ws.Range("A5").Value = "Stiven"
$ws.Range("C5").Value = "ADMIN"
